$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 39
$ws.Range("E3").Value = 18
$ws.Range("E17").Value = 89
$ws.Range("E24").Value = 19
$ws.Range("E55").Value = 6
$ws.Range("E60").Value = 14
$ws.Range("E62").Value = 35
$ws.Range("F62").Value = 7
$ws.Range("H62").Value = 7
$ws.Range("E63").Value = 22
$ws.Range("E70").Value = 31
$ws.Range("E76").Value = 37
$ws.Range("E89").Value = 24
